$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-14 Saturday" "2023-10-15 Sunday"

Replace-Text "23×36=" "97×53="
Replace-Text "73×91=" "25×99="
Replace-Text "37×12=" "49×99="
Replace-Text "79×11=" "46×32="
Replace-Text "87×27=" "95×34="

Replace-Text "58×89=" "51×88="
Replace-Text "80×33=" "40×48="
Replace-Text "91×11=" "75×78="
Replace-Text "96×37=" "77×18="
Replace-Text "69×46=" "92×79="

Replace-Text "59×73=" "30×57="
Replace-Text "22×54=" "60×72="
Replace-Text "12×96=" "88×18="
Replace-Text "13×68=" "43×36="
Replace-Text "62×32=" "24×36="

Replace-Text "11×42=" "65×19="
Replace-Text "71×72=" "84×86="
Replace-Text "23×67=" "30×14="
Replace-Text "73×61=" "63×28="
Replace-Text "21×26=" "96×34="

Replace-Text "39×21=" "74×50="
Replace-Text "52×51=" "56×63="
Replace-Text "67×56=" "19×37="
Replace-Text "81×36=" "56×82="
Replace-Text "43×15=" "81×76="
